# The edit removes the (unused) built-in "Footnote Text" paragraph style
# definition from styles.xml. It was present in the source document but is
# not referenced by any paragraph, so the rebuilt document simply drops the
# explicit style definition (Word still recognises the built-in name, it
# just no longer carries an explicit <w:style> entry in styles.xml).
$d = $word.ActiveDocument

$footnoteStyle = $d.Styles("Footnote Text")
try {
    $footnoteStyle.Delete()
} catch {
    # Already absent from the document's explicit style collection -
    # nothing further to do.
}
